# Insert a new weekly record at row 10 (pushing the existing rows 10-47 down
# to 11-48). The new row duplicates the surrounding record's descriptive
# columns (market / region / product / quality / unit / origin / etc.) and
# only carries a new date + price figures for this extra reporting week.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 10:47 down to 11:48, leaving a blank row 10 to populate.
$ws.Rows("10:10").Insert()

$newRow = 10

$ws.Cells.Item($newRow, 1).Value  = 1
$ws.Cells.Item($newRow, 2).Value  = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item($newRow, 3).Value  = "Arica y Parinacota"
$ws.Cells.Item($newRow, 4).Value  = 44525
$ws.Cells.Item($newRow, 5).Value  = 15
$ws.Cells.Item($newRow, 6).Value  = 100112012
$ws.Cells.Item($newRow, 7).Value  = "Espinaca"
$ws.Cells.Item($newRow, 8).Value  = "Sin especificar"
$ws.Cells.Item($newRow, 9).Value  = "Primera"
$ws.Cells.Item($newRow, 10).Value = 250
$ws.Cells.Item($newRow, 11).Value = 1800
$ws.Cells.Item($newRow, 12).Value = 2000
$ws.Cells.Item($newRow, 13).Value = 1900
$ws.Cells.Item($newRow, 14).Value = "$/atado 2,5 a 3 kilos"
$ws.Cells.Item($newRow, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item($newRow, 16).Value = 633
$ws.Cells.Item($newRow, 17).Value = 3
$ws.Cells.Item($newRow, 18).Value = "Hortaliza"
